# Apply the Welcome.docx (Russian) content edits.
$d = $word.ActiveDocument

# 1. Trim the trailing clause from the ASIC/mining paragraph.
$d.Content.Find.Execute(
    ", until Smartcash reaches a considerable market cap.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, ".", 2
) | Out-Null

# 2. Recreate the "exchanges" bookmark so it is re-registered
#    (name / anchor position are preserved; Word renumbers the
#    underlying bookmark id bookkeeping on save).
$bm = $d.Bookmarks("exchanges")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("exchanges", $bmRange) | Out-Null
